$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the empty heading paragraph
#    down to right after "...Science, 2nd " (before "student "),
#    i.e. right before the word "student". Do this first so the
#    later text edits in this same paragraph don't re-merge runs
#    across the insertion point.
# ------------------------------------------------------------------
$rStudent = $d.Content
$foundStudent = $rStudent.Find.Execute("student ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($foundStudent) {
    $rStudent.Collapse(1)   # wdCollapseStart
    $d.Bookmarks.Add("_GoBack", $rStudent)
}

# ------------------------------------------------------------------
# 2) "Studying Honor Bachelor of Science, 2[nd]" -> "...1[st]"
#    (2nd year -> 1st year). Narrowly scoped so no other "nd"/"2"
#    occurrence elsewhere in the document is touched.
# ------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Studying Honor Bachelor of Science, 2", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    # last character of the match is the "2" -> change to "1"
    $rDigit = $d.Range($r.End - 1, $r.End)
    $rDigit.Text = "1"

    # the following two characters are the superscript "nd" -> "st"
    $rSuffix = $d.Range($r.End, $r.End + 2)
    $rSuffix.Text = "st"
}

# ------------------------------------------------------------------
# 3) "Current GPD: 4.0 /" -> "Cumulative GPA: 3.96 /"
# ------------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Current GPD: 4.0 /", $true, $false, $false, $false, $false, $true, 1, $false, "Cumulative GPA: 3.96 /", 2) | Out-Null
